$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.236.21'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.27%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.855.08'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.16%  '

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7011'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +2.58%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '237.98'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.30%  '

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.08043'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +4.90%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3024'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.36%  '

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +2.13%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08178'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.53%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.868.28'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.63%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.202'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.50%  '

$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.7069'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -2.07%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '89.73'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.56%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '29.099.88'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -0.14%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.829'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +2.10%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000007857'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.74%  '

$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.83%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '236.35'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +1.28%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -0.04%  '

$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.06%  '

$ws.Range('B23').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.042.16'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -2.69%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '7.499'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.30%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '163.31'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +0.99%  '

$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.62%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1407'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -1.29%  '

$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.41%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.913'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -2.04%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.408'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +1.08%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.471'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.86%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.360'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -3.29%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.026'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +0.58%  '

$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.83%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.164'

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7150'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +1.87%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9961'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -2.81%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.685'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.42%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01846'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.35%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.722'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +1.69%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9342'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +2.90%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.153.03'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +4.68%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '6.003'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +0.56%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.4261'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -0.13%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '70.27'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.64%  '

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.05%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '102.80'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +0.72%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5284'

$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -1.12%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.995.69'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +0.31%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '9.149'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.19%  '
